# Auto-generated Excel COM script to apply value updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 784.75
$ws.Cells.Item(43, 10).Value = 718.4
$ws.Cells.Item(43, 12).Value = 718.4
$ws.Cells.Item(43, 14).Value = -856.4
$ws.Cells.Item(86, 8).Value = 4239
$ws.Cells.Item(86, 9).Value = 4337.9
$ws.Cells.Item(86, 10).Value = 3250
$ws.Cells.Item(86, 11).Value = 4337.9
$ws.Cells.Item(86, 12).Value = 3250
$ws.Cells.Item(86, 13).Value = -3214.9
$ws.Cells.Item(86, 14).Value = -5496
$ws.Cells.Item(88, 8).Value = 25799.4
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 25799.4
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 25799.4
$ws.Cells.Item(88, 13).Value = $null
$ws.Cells.Item(88, 14).Value = -26611.4
$ws.Cells.Item(89, 8).Value = 4239
$ws.Cells.Item(89, 9).Value = 4337.9
$ws.Cells.Item(89, 10).Value = 3250
$ws.Cells.Item(89, 11).Value = 21689.5
$ws.Cells.Item(89, 12).Value = 16250
$ws.Cells.Item(89, 13).Value = -16073.5
$ws.Cells.Item(89, 14).Value = -27482
$ws.Cells.Item(91, 8).Value = 25799.4
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 25799.4
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 25799.4
$ws.Cells.Item(91, 13).Value = $null
$ws.Cells.Item(91, 14).Value = -28607.4
$ws.Cells.Item(134, 8).Value = 56573.332
$ws.Cells.Item(134, 10).Value = 56573.332
$ws.Cells.Item(134, 12).Value = 56573.332
$ws.Cells.Item(134, 14).Value = -66713.33199999999
$ws.Cells.Item(137, 8).Value = 1458.5555
$ws.Cells.Item(137, 9).Value = 1514.1578
$ws.Cells.Item(137, 10).Value = 1326.5
$ws.Cells.Item(137, 11).Value = 4542.4734
$ws.Cells.Item(137, 12).Value = 3979.5
$ws.Cells.Item(137, 13).Value = -1992.4734
$ws.Cells.Item(137, 14).Value = -9079.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 12835512
$ws.Cells.Item(122, 9).Value = 12835512
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 38506536
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -38504086
$ws.Cells.Item(122, 14).Value = $null
$ws.Cells.Item(135, 8).Value = 59879.9
$ws.Cells.Item(135, 10).Value = 59879.9
$ws.Cells.Item(135, 12).Value = 59879.9
$ws.Cells.Item(135, 14).Value = -70019.89999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(56, 8).Value = 19110
$ws.Cells.Item(56, 10).Value = 19110
$ws.Cells.Item(56, 12).Value = 19110
$ws.Cells.Item(56, 14).Value = -20588
$ws.Cells.Item(134, 8).Value = 33350.25
$ws.Cells.Item(134, 9).Value = 2109.8518
$ws.Cells.Item(134, 10).Value = 202048.4
$ws.Cells.Item(134, 11).Value = 6329.555399999999
$ws.Cells.Item(134, 12).Value = 606145.2
$ws.Cells.Item(134, 13).Value = -3794.555399999999
$ws.Cells.Item(134, 14).Value = -611215.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2555.889
$ws.Cells.Item(16, 9).Value = 2460.476
$ws.Cells.Item(16, 11).Value = 2460.476
$ws.Cells.Item(16, 13).Value = -2173.476
$ws.Cells.Item(76, 8).Value = 335693
$ws.Cells.Item(76, 9).Value = 335693
$ws.Cells.Item(76, 11).Value = 335693
$ws.Cells.Item(76, 13).Value = -335378
$ws.Cells.Item(79, 8).Value = 335693
$ws.Cells.Item(79, 9).Value = 335693
$ws.Cells.Item(79, 11).Value = 335693
$ws.Cells.Item(79, 13).Value = -334601
$ws.Cells.Item(99, 8).Value = 5000
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 5000
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 5000
$ws.Cells.Item(99, 13).Value = $null
$ws.Cells.Item(99, 14).Value = -7996
$ws.Cells.Item(113, 8).Value = 2555.889
$ws.Cells.Item(113, 9).Value = 2460.476
$ws.Cells.Item(113, 11).Value = 2460.476
$ws.Cells.Item(113, 13).Value = -290.4760000000001
$ws.Cells.Item(126, 8).Value = 5000
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 13).Value = $null
$ws.Cells.Item(126, 14).Value = -19940
$ws.Cells.Item(135, 8).Value = 38682
$ws.Cells.Item(135, 10).Value = 38682
$ws.Cells.Item(135, 12).Value = 38682
$ws.Cells.Item(135, 14).Value = -48822

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 164425.67
$ws.Cells.Item(113, 9).Value = 479.875
$ws.Cells.Item(113, 10).Value = 270768.9
$ws.Cells.Item(113, 11).Value = 1439.625
$ws.Cells.Item(113, 12).Value = 812306.7000000001
$ws.Cells.Item(113, 13).Value = 730.375
$ws.Cells.Item(113, 14).Value = -816646.7000000001
$ws.Cells.Item(122, 8).Value = 6140.619
$ws.Cells.Item(122, 9).Value = 1066
$ws.Cells.Item(122, 10).Value = 11722.7
$ws.Cells.Item(122, 11).Value = 9594
$ws.Cells.Item(122, 12).Value = 105504.3
$ws.Cells.Item(122, 13).Value = -7144
$ws.Cells.Item(122, 14).Value = -110404.3
$ws.Cells.Item(131, 8).Value = 2174933.5
$ws.Cells.Item(131, 9).Value = 5883005
$ws.Cells.Item(131, 10).Value = 1236.5862
$ws.Cells.Item(131, 11).Value = 17649015
$ws.Cells.Item(131, 12).Value = 3709.7586
$ws.Cells.Item(131, 13).Value = -17643975
$ws.Cells.Item(131, 14).Value = -13789.7586

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 9515
$ws.Cells.Item(80, 9).Value = 10642.5
$ws.Cells.Item(80, 11).Value = 10642.5
$ws.Cells.Item(80, 13).Value = -9644.5
$ws.Cells.Item(83, 8).Value = 9515
$ws.Cells.Item(83, 9).Value = 10642.5
$ws.Cells.Item(83, 11).Value = 53212.5
$ws.Cells.Item(83, 13).Value = -48220.5
$ws.Cells.Item(113, 8).Value = 1000000000
$ws.Cells.Item(113, 9).Value = 1000000000
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 1000000000
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -999997830
$ws.Cells.Item(113, 14).Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 6071
$ws.Cells.Item(32, 9).Value = 6071
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 6071
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -5754
$ws.Cells.Item(32, 14).Value = $null
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).Value = $null
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = $null
$ws.Cells.Item(58, 14).Value = $null
$ws.Cells.Item(61, 8).Value = 4000
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 4000
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 4000
$ws.Cells.Item(61, 13).Value = $null
$ws.Cells.Item(61, 14).Value = -4404
$ws.Cells.Item(68, 8).Value = 38463604
$ws.Cells.Item(68, 9).Value = 1952.3077
$ws.Cells.Item(68, 10).Value = 76925256
$ws.Cells.Item(68, 11).Value = 1952.3077
$ws.Cells.Item(68, 12).Value = 76925256
$ws.Cells.Item(68, 13).Value = -1203.3077
$ws.Cells.Item(68, 14).Value = -76926754
$ws.Cells.Item(71, 8).Value = 38463604
$ws.Cells.Item(71, 9).Value = 1952.3077
$ws.Cells.Item(71, 10).Value = 76925256
$ws.Cells.Item(71, 11).Value = 9761.538500000001
$ws.Cells.Item(71, 12).Value = 384626280
$ws.Cells.Item(71, 13).Value = -6017.538500000001
$ws.Cells.Item(71, 14).Value = -384633768
$ws.Cells.Item(82, 8).Value = 556614
$ws.Cells.Item(82, 9).Value = 910597.8
$ws.Cells.Item(82, 10).Value = 123967.11
$ws.Cells.Item(82, 11).Value = 910597.8
$ws.Cells.Item(82, 12).Value = 123967.11
$ws.Cells.Item(82, 13).Value = -910236.8
$ws.Cells.Item(82, 14).Value = -124689.11
$ws.Cells.Item(85, 8).Value = 556614
$ws.Cells.Item(85, 9).Value = 910597.8
$ws.Cells.Item(85, 10).Value = 123967.11
$ws.Cells.Item(85, 11).Value = 910597.8
$ws.Cells.Item(85, 12).Value = 123967.11
$ws.Cells.Item(85, 13).Value = -909349.8
$ws.Cells.Item(85, 14).Value = -126463.11
$ws.Cells.Item(100, 8).Value = 2422.1365
$ws.Cells.Item(100, 9).Value = 2267.1667
$ws.Cells.Item(100, 10).Value = 2480.25
$ws.Cells.Item(100, 11).Value = 2267.1667
$ws.Cells.Item(100, 12).Value = 2480.25
$ws.Cells.Item(100, 13).Value = -1726.1667
$ws.Cells.Item(100, 14).Value = -3562.25
$ws.Cells.Item(113, 8).Value = 4000
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 4000
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 4000
$ws.Cells.Item(113, 13).Value = $null
$ws.Cells.Item(113, 14).Value = -8340
$ws.Cells.Item(122, 8).Value = 40714284
$ws.Cells.Item(122, 9).Value = 71428570
$ws.Cells.Item(122, 10).Value = 10000000
$ws.Cells.Item(122, 11).Value = 214285710
$ws.Cells.Item(122, 12).Value = 30000000
$ws.Cells.Item(122, 13).Value = -214283260
$ws.Cells.Item(122, 14).Value = -30004900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 36365.223
$ws.Cells.Item(46, 10).Value = 36365.223
$ws.Cells.Item(46, 12).Value = 36365.223
$ws.Cells.Item(46, 14).Value = -36827.223
$ws.Cells.Item(92, 8).Value = 42550
$ws.Cells.Item(92, 10).Value = 42550
$ws.Cells.Item(92, 12).Value = 42550
$ws.Cells.Item(92, 14).Value = -47542
$ws.Cells.Item(113, 8).Value = 1533
$ws.Cells.Item(113, 9).Value = 1533
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 4599
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -2429
$ws.Cells.Item(113, 14).Value = $null
$ws.Cells.Item(134, 8).Value = 36365.223
$ws.Cells.Item(134, 10).Value = 36365.223
$ws.Cells.Item(134, 12).Value = 109095.669
$ws.Cells.Item(134, 14).Value = -114165.669
$ws.Cells.Item(135, 8).Value = 47339.918
$ws.Cells.Item(135, 10).Value = 47339.918
$ws.Cells.Item(135, 12).Value = 47339.918
$ws.Cells.Item(135, 14).Value = -57479.918
